$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.903.68'
$ws.Range("D2").ClearFormats()

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.70%  '
$ws.Range("E2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.683.51'
$ws.Range("D3").ClearFormats()

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("E4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.50'
$ws.Range("D5").ClearFormats()

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("E5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.85'
$ws.Range("D6").ClearFormats()

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.58%  '
$ws.Range("E6").ClearFormats()

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E7").ClearFormats()

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.596'
$ws.Range("D8").ClearFormats()

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.20%  '
$ws.Range("E8").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.96'
$ws.Range("D10").ClearFormats()

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.61%  '
$ws.Range("E10").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.399'
$ws.Range("D11").ClearFormats()

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.18%  '
$ws.Range("E11").ClearFormats()

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("E12").ClearFormats()

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.54'
$ws.Range("D13").ClearFormats()

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.29%  '
$ws.Range("E13").ClearFormats()

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000201'
$ws.Range("D14").ClearFormats()

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +5.68%  '
$ws.Range("E14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.166.19'
$ws.Range("D15").ClearFormats()

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.94%  '
$ws.Range("E15").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.691.23'
$ws.Range("D16").ClearFormats()

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("E16").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.673.88'
$ws.Range("D17").ClearFormats()

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("E17").ClearFormats()

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.63'
$ws.Range("D18").ClearFormats()

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.11%  '
$ws.Range("E18").ClearFormats()

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.07%  '
$ws.Range("E19").ClearFormats()

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.20%  '
$ws.Range("E20").ClearFormats()

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '352.83'
$ws.Range("D21").ClearFormats()

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.99%  '
$ws.Range("E21").ClearFormats()

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("E22").ClearFormats()

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.53'
$ws.Range("D23").ClearFormats()

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.40%  '
$ws.Range("E23").ClearFormats()

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +7.11%  '
$ws.Range("E24").ClearFormats()

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.87'
$ws.Range("D25").ClearFormats()

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +3.64%  '
$ws.Range("E25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.64'
$ws.Range("D26").ClearFormats()

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.44%  '
$ws.Range("E26").ClearFormats()

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.43%  '
$ws.Range("E27").ClearFormats()

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.31%  '
$ws.Range("E28").ClearFormats()

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.97%  '
$ws.Range("E29").ClearFormats()

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.14%  '
$ws.Range("E30").ClearFormats()

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.79%  '
$ws.Range("E31").ClearFormats()

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '536.09'
$ws.Range("D32").ClearFormats()

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.63%  '
$ws.Range("E32").ClearFormats()

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.78'
$ws.Range("D33").ClearFormats()

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.39%  '
$ws.Range("E33").ClearFormats()

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.59'
$ws.Range("D34").ClearFormats()

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.54%  '
$ws.Range("E34").ClearFormats()

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.25%  '
$ws.Range("E35").ClearFormats()

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.62%  '
$ws.Range("E36").ClearFormats()

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.50%  '
$ws.Range("E37").ClearFormats()

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '160.56'
$ws.Range("D38").ClearFormats()

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.49%  '
$ws.Range("E38").ClearFormats()

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("E39").ClearFormats()

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.00%  '
$ws.Range("E40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").ClearFormats()

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("E41").ClearFormats()

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.34'
$ws.Range("D42").ClearFormats()

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("E42").ClearFormats()

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '166.67'
$ws.Range("D43").ClearFormats()

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("E43").ClearFormats()

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.12'
$ws.Range("D44").ClearFormats()

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.72%  '
$ws.Range("E44").ClearFormats()

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0621'
$ws.Range("D45").ClearFormats()

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("E45").ClearFormats()

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.22'
$ws.Range("D46").ClearFormats()

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.22%  '
$ws.Range("E46").ClearFormats()

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.89%  '
$ws.Range("E47").ClearFormats()

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("E48").ClearFormats()

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.652'
$ws.Range("D49").ClearFormats()

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.28%  '
$ws.Range("E49").ClearFormats()

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.40'
$ws.Range("D50").ClearFormats()

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.99%  '
$ws.Range("E50").ClearFormats()

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0996'
$ws.Range("D51").ClearFormats()

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.90%  '
$ws.Range("E51").ClearFormats()
